$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column C for the new rows as text so date-like strings are not
# auto-converted to date serials, then revert the style so no stray
# cell formatting is left behind (matches the source file which has no
# explicit styles on these cells).
$dateRange = $ws.Range("C94:C126")
$dateRange.NumberFormat = "@"

$ws.Range("C94").Value = "2019-11-01"
$ws.Range("D94").Value = 228
$ws.Range("E94").Value = 296
$ws.Range("F94").Value = 183
$ws.Range("G94").Value = 284
$ws.Range("H94").Value = 371

$ws.Range("C95").Value = "2019-11-02"
$ws.Range("D95").Value = 228
$ws.Range("E95").Value = 296
$ws.Range("F95").Value = 178
$ws.Range("G95").Value = 282
$ws.Range("H95").Value = 366

$ws.Range("C96").Value = "2019-11-01"
$ws.Range("D96").Value = 228
$ws.Range("E96").Value = 296

$ws.Range("C97").Value = "2019-11-02"

$ws.Range("C98").Value = "2019-11-01"
$ws.Range("D98").Value = 228
$ws.Range("E98").Value = 296

$ws.Range("C99").Value = "2019-11-02"

$ws.Range("C100").Value = "2019-11-01"
$ws.Range("D100").Value = 228
$ws.Range("E100").Value = 296

$ws.Range("C101").Value = "2019-11-02"

$ws.Range("C102").Value = "2019-11-01"
$ws.Range("D102").Value = 228
$ws.Range("E102").Value = 296

$ws.Range("C103").Value = "2019-11-02"

$ws.Range("C104").Value = "2019-07-15"
$ws.Range("F104").Value = 321

$ws.Range("C105").Value = "2019-07-18"
$ws.Range("F105").Value = 321

$ws.Range("C106").Value = "2019-07-21"
$ws.Range("F106").Value = 321

$ws.Range("C107").Value = "2019-07-24"
$ws.Range("F107").Value = 321

$ws.Range("C108").Value = "2019-07-27"

$ws.Range("C109").Value = "2019-07-30"

$ws.Range("C110").Value = "2019-08-02"

$ws.Range("C111").Value = "2019-08-05"

$ws.Range("C112").Value = "2019-08-08"

$ws.Range("C113").Value = "2019-08-11"
$ws.Range("F113").Value = 321

$ws.Range("C114").Value = "2019-08-14"
$ws.Range("F114").Value = 321

$ws.Range("C115").Value = "2019-08-15"

$ws.Range("C116").Value = "2019-08-01"
$ws.Range("D116").Value = 410

$ws.Range("C117").Value = "2019-08-04"
$ws.Range("D117").Value = 410

$ws.Range("C118").Value = "2019-08-07"
$ws.Range("D118").Value = 410

$ws.Range("C119").Value = "2019-08-10"
$ws.Range("D119").Value = 410

$ws.Range("C120").Value = "2019-08-13"
$ws.Range("D120").Value = 410
$ws.Range("F120").Value = 330

$ws.Range("C121").Value = "2019-08-16"
$ws.Range("D121").Value = 410
$ws.Range("F121").Value = 330

$ws.Range("C122").Value = "2019-08-19"
$ws.Range("D122").Value = 410
$ws.Range("F122").Value = 330

$ws.Range("C123").Value = "2019-08-22"

$ws.Range("C124").Value = "2019-08-25"
$ws.Range("F124").Value = 331

$ws.Range("C125").Value = "2019-08-28"
$ws.Range("F125").Value = 331

$ws.Range("C126").Value = "2019-08-31"

$dateRange.Style = "Normal"